$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header cell in column H, matching the formatting of the other
# header cells (bold / bordered / centered), same as the rest of row 1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cells (rows 2-7) for the "Save" column, all zero, using the
# default (unstyled) formatting like the rest of the numeric columns.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
